# UserStories.xlsx — spelling/wording fixes on the "SprintBacklog2" sheet
# plus a duplicate "page d'accueil" user-story merged into the earlier
# (correctly spelled) shared string, and a row-height / selection touch-up.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SprintBacklog2")

# --- Text corrections (accents, typos, wording) -----------------------
$ws.Range("B4").Value  = "ajout dans la logique promesse don(méthodes)"
$ws.Range("B5").Value  = "form validation(étoiles, message d'erreur, formulaire vide)"
$ws.Range("B6").Value  = "modal redirection après promesse don reçu"
$ws.Range("B7").Value  = "mode de livraison interface, modification logique pour offrir don(2 radio butons et leur désactivation si argent)"
$ws.Range("B9").Value  = "modification dans la table concernant l'adresse"
$ws.Range("B10").Value = "modification logique pour offrir don pour les champ ajouter"
$ws.Range("B12").Value = "interface - liste des dons qui sont pas traités(nom, quantité, catégorie, montant, mode livraison, date promise)"

# Row 17 used a misspelled duplicate of "page d'accueil" (already present
# earlier in the shared-string table, e.g. row 3 of SprintBacklog_1). Reuse
# the correctly spelled text so the duplicate string collapses away.
$ws.Range("A17").Value = "page d'accueil"
$ws.Range("B17").Value = "modification – button, caroussele ?"

$ws.Range("B21").Value = "envoie courriel a partir de lien modifier"
$ws.Range("B22").Value = "se authentifier"
$ws.Range("A23").Value = "propose se disponibilité comme bénévole"
$ws.Range("B23").Value = "button nav pour formulaire"
$ws.Range("B26").Value = "classe, DAO"

# --- Row heights (re-measured after the text edits above) -------------
$ws.Rows.Item(7).RowHeight = 23.85
$ws.Rows.Item(12).RowHeight = 23.85

# --- View state: scroll back to the top and leave the selection on the
#     last row, and make sure this sheet stays the active tab ----------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B27").Select()
